$d = $word.ActiveDocument

# Step 1: append two trailing spaces to the first paragraph's existing text
# by replacing it in place (keeps it as a single run).
$rng = $d.Content
$found = $rng.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

# Step 2: insert a new run right after the replaced text (before the
# paragraph mark) containing the red "change" annotation.
$insertRange = $d.Range($rng.End, $rng.End)
$insertRange.InsertAfter("(This is a change " + [char]8211 + " Version for branch alternate)")
$insertRange.Font.Color = 192
